# This sheet is a weekly price log for "Cebollín baby" (Agrícola del Norte
# S.A. de Arica). A new weekly observation was added at the top of the data
# block (row 21): every existing row from 21..53 shifts its date/volume/price
# figures down into the next row (row N <- old row N-1), row 54 is a brand
# new row holding what used to be row 53's data, and row 21 itself receives
# the genuinely new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D21").Value = 44498
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 900
$ws.Range("P21").Value = 450

$ws.Range("D22").Value = 44417
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3500
$ws.Range("M22").Value = 3250
$ws.Range("P22").Value = 1625

$ws.Range("D23").Value = 44424
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 2500
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = 2750
$ws.Range("P23").Value = 1375

$ws.Range("D24").Value = 44495
$ws.Range("J24").Value = 300
$ws.Range("K24").Value = 1000
$ws.Range("L24").Value = 1200
$ws.Range("M24").Value = 1100
$ws.Range("P24").Value = 550

$ws.Range("D25").Value = 44370
$ws.Range("J25").Value = 400
$ws.Range("K25").Value = 3400
$ws.Range("L25").Value = 3500
$ws.Range("M25").Value = 3445
$ws.Range("P25").Value = 1722

$ws.Range("D26").Value = 44284
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 1800
$ws.Range("L26").Value = 2000
$ws.Range("M26").Value = 1900
$ws.Range("P26").Value = 950

$ws.Range("D27").Value = 44428
$ws.Range("J27").Value = 270
$ws.Range("K27").Value = 3500
$ws.Range("L27").Value = 3800
$ws.Range("M27").Value = 3650
$ws.Range("P27").Value = 1825

$ws.Range("D28").Value = 44293
$ws.Range("J28").Value = 250
$ws.Range("K28").Value = 1500
$ws.Range("L28").Value = 1800
$ws.Range("M28").Value = 1650
$ws.Range("P28").Value = 825

$ws.Range("D29").Value = 44343
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 1500
$ws.Range("L29").Value = 2000
$ws.Range("M29").Value = 1750
$ws.Range("P29").Value = 875

$ws.Range("D30").Value = 44484
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 950
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = 975
$ws.Range("P30").Value = 488

$ws.Range("D31").Value = 44356
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 2400
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = 2450
$ws.Range("P31").Value = 1225

$ws.Range("D32").Value = 44169
$ws.Range("J32").Value = 300
$ws.Range("K32").Value = 2000
$ws.Range("L32").Value = 2500
$ws.Range("M32").Value = 2250
$ws.Range("P32").Value = 1125

$ws.Range("D33").Value = 44410
$ws.Range("J33").Value = 250
$ws.Range("K33").Value = 2800
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 2900
$ws.Range("P33").Value = 1450

$ws.Range("D34").Value = 44319
$ws.Range("J34").Value = 300
$ws.Range("K34").Value = 1900
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 1950
$ws.Range("P34").Value = 975

$ws.Range("D35").Value = 44473
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 950
$ws.Range("L35").Value = 1000
$ws.Range("M35").Value = 975
$ws.Range("P35").Value = 488

$ws.Range("D36").Value = 44235
$ws.Range("J36").Value = 250
$ws.Range("K36").Value = 4500
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = 4750
$ws.Range("P36").Value = 2375

$ws.Range("D37").Value = 44203
$ws.Range("J37").Value = 300
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 2500
$ws.Range("M37").Value = 2250
$ws.Range("P37").Value = 1125

$ws.Range("D38").Value = 44469
$ws.Range("J38").Value = 250
$ws.Range("K38").Value = 900
$ws.Range("L38").Value = 1000
$ws.Range("M38").Value = 950
$ws.Range("P38").Value = 475

$ws.Range("D39").Value = 44165
$ws.Range("J39").Value = 300
$ws.Range("K39").Value = 1000
$ws.Range("L39").Value = 1200
$ws.Range("M39").Value = 1100
$ws.Range("P39").Value = 550

$ws.Range("D40").Value = 44320
$ws.Range("J40").Value = 250
$ws.Range("K40").Value = 1400
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = 1450
$ws.Range("P40").Value = 725

$ws.Range("D41").Value = 44258
$ws.Range("J41").Value = 150
$ws.Range("K41").Value = 2400
$ws.Range("L41").Value = 2500
$ws.Range("M41").Value = 2450
$ws.Range("P41").Value = 1225

$ws.Range("D42").Value = 44349
$ws.Range("J42").Value = 300
$ws.Range("K42").Value = 1800
$ws.Range("L42").Value = 2000
$ws.Range("M42").Value = 1900
$ws.Range("P42").Value = 950

$ws.Range("D43").Value = 44448
$ws.Range("J43").Value = 270
$ws.Range("K43").Value = 1900
$ws.Range("L43").Value = 2000
$ws.Range("M43").Value = 1950
$ws.Range("P43").Value = 975

$ws.Range("D44").Value = 44243
$ws.Range("J44").Value = 200
$ws.Range("K44").Value = 2900
$ws.Range("L44").Value = 3000
$ws.Range("M44").Value = 2950
$ws.Range("P44").Value = 1475

$ws.Range("D45").Value = 44452
$ws.Range("J45").Value = 300
$ws.Range("K45").Value = 1900
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = 1950
$ws.Range("P45").Value = 975

$ws.Range("D46").Value = 44435
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 1800
$ws.Range("L46").Value = 2000
$ws.Range("M46").Value = 1930
$ws.Range("P46").Value = 965

$ws.Range("D47").Value = 44433
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 1800
$ws.Range("L47").Value = 2000
$ws.Range("M47").Value = 1900
$ws.Range("P47").Value = 950

$ws.Range("D48").Value = 44334
$ws.Range("J48").Value = 200
$ws.Range("K48").Value = 2800
$ws.Range("L48").Value = 3000
$ws.Range("M48").Value = 2900
$ws.Range("P48").Value = 1450

$ws.Range("D49").Value = 44386
$ws.Range("J49").Value = 250
$ws.Range("K49").Value = 3500
$ws.Range("L49").Value = 4000
$ws.Range("M49").Value = 3750
$ws.Range("P49").Value = 1875

$ws.Range("D50").Value = 44263
$ws.Range("J50").Value = 270
$ws.Range("K50").Value = 1900
$ws.Range("L50").Value = 2000
$ws.Range("M50").Value = 1950
$ws.Range("P50").Value = 975

$ws.Range("D51").Value = 44298
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 1400
$ws.Range("L51").Value = 1500
$ws.Range("M51").Value = 1450
$ws.Range("P51").Value = 725

$ws.Range("D52").Value = 44397
$ws.Range("J52").Value = 300
$ws.Range("K52").Value = 3500
$ws.Range("L52").Value = 4000
$ws.Range("M52").Value = 3750
$ws.Range("P52").Value = 1875

$ws.Range("D53").Value = 44312
$ws.Range("J53").Value = 300
$ws.Range("K53").Value = 1000
$ws.Range("L53").Value = 1200
$ws.Range("M53").Value = 1100
$ws.Range("P53").Value = 550

# Row 54 is brand new (the dimension grows from A1:R53 to A1:R54). It carries
# the same constant descriptive fields as every other data row plus the
# price figures that used to live in row 53.
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44326
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100112038
$ws.Range("G54").Value = "Cebollín baby"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 300
$ws.Range("K54").Value = 1400
$ws.Range("L54").Value = 1500
$ws.Range("M54").Value = 1450
$ws.Range("N54").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 725
$ws.Range("Q54").Value = 2
$ws.Range("R54").Value = "Hortaliza"
